$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$values = @{
    2 = 34
    4 = 26
    5 = 1132
    6 = 8756
    8 = 231
    9 = 7078
    10 = 183
    11 = 312
    12 = 5331
    13 = 68
    14 = 5991
    15 = 1090
    16 = 383
    17 = 392
    18 = 547
    19 = 323
    20 = 266
    21 = 142
    22 = 200
    23 = 160
    25 = 9858
    26 = 79
    27 = 1852
    28 = 1645
    29 = 47
    30 = 36
    31 = 1992
    32 = 77
    33 = 83
    34 = 1025
    35 = 2011
    36 = 300
    37 = 1335
    38 = 58
    39 = 5117
    40 = 1191
    41 = 624
    43 = 167
    44 = 1095
    45 = 1067
    46 = 964
    47 = 1338
    48 = 58
    49 = 1089
}
foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}

$ws = $wb.Worksheets.Item("演出")
$values = @{
    2 = 37
    4 = 14
    5 = 5
    10 = 190
    12 = 9
    18 = 902
}
foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}

$ws = $wb.Worksheets.Item("全部类型")
$values = @{
    4 = 63
    5 = 26
    6 = 14
    7 = 1132
    8 = 8756
    10 = 231
    11 = 7078
    12 = 183
    13 = 312
    16 = 5331
    17 = 5331
    18 = 5991
    19 = 5991
    20 = 1090
    21 = 383
    22 = 392
    23 = 547
    24 = 323
    25 = 266
    28 = 160
    29 = 190
    30 = 9858
    32 = 1852
    33 = 1645
    34 = 47
    35 = 1992
    36 = 77
    37 = 83
    38 = 1025
    39 = 2011
    40 = 300
    41 = 1335
    42 = 5117
    43 = 1191
    44 = 624
    45 = 99
    46 = 167
    47 = 1095
    48 = 1067
    49 = 964
    50 = 1338
    51 = 1089
}
foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}

